# "Fruta / hortaliza, semanal" weekly data refresh.
#
# A new weekly price record is inserted into the Cilantro price table at
# row 299, pushing the existing rows 299-376 down to 300-377 (the sheet
# grows from 376 to 377 rows; dimension A1:R376 -> A1:R377).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 299 - Excel's default
# Insert() shifts rows 299:376 down to 300:377 and the new row inherits
# row 298's formatting (so column D keeps its date number format).
$ws.Rows.Item(299).Insert()

# Populate the newly inserted row 299 with the new record's data.
# (Single-quoted literals throughout so none of the text - e.g. the
# leading "$/" in the unit string - is ever treated as interpolation.)
$ws.Range("A299").Value = 10
$ws.Range("B299").Value = 'Vega Modelo de Temuco'
$ws.Range("C299").Value = 'La Araucanía'
$ws.Range("D299").Value = 44754
$ws.Range("E299").Value = 9
$ws.Range("F299").Value = 100112040
$ws.Range("G299").Value = 'Cilantro'
$ws.Range("H299").Value = 'Sin especificar'
$ws.Range("I299").Value = 'Primera'
$ws.Range("J299").Value = 50
$ws.Range("K299").Value = 4300
$ws.Range("L299").Value = 4300
$ws.Range("M299").Value = 4300
$ws.Range("N299").Value = '$/docena de atados (2 kilos)'
$ws.Range("O299").Value = 'Provincia de Cautín'
$ws.Range("P299").Value = 2150
$ws.Range("Q299").Value = 2
$ws.Range("R299").Value = 'Hortaliza'
